$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data: a new "lastname" (shared string "Matin") and its
# corresponding "chat_id" value, mirroring the existing username/lastname/
# chat_id rows already on the sheet.
$ws.Range("C3").Value = "Matin"
$ws.Range("D3").Value = 8069824403
